# Update Name of Algo - apply refreshed RandomForest imputation results
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "B4"   = 8.7041
    "B7"   = 4.830999999999997
    "C7"   = -13.74569999999999
    "C15"  = -14.666
    "B16"  = 6.794099999999997
    "C21"  = -12.16830000000001
    "C22"  = -12.66170000000001
    "C23"  = -12.37810000000001
    "B28"  = 5.634500000000003
    "B29"  = 5.128600000000002
    "B32"  = 7.186099999999997
    "C34"  = -11.61260000000002
    "B40"  = 8.728799999999993
    "C43"  = -12.8817
    "C45"  = -13.76729999999999
    "C50"  = -13.94889999999999
    "C51"  = -11.79589999999999
    "B52"  = 5.231199999999999
    "B57"  = 4.870099999999995
    "B66"  = 5.8104
    "C66"  = -11.6921
    "C67"  = -11.5865
    "C79"  = -11.45920000000001
    "C84"  = -13.14809999999999
    "C92"  = -11.3408
    "C97"  = -12.5955
    "B100" = 5.415099999999996
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
